# Remove the yellow highlighting that was applied to the
# "SlotDiffusion: Object-Centric Generative Modeling with Diffusion Models"
# paper title (the run(s) carrying <w:highlight w:val="yellow"/> in
# <w:rPr/>). Other highlighted text elsewhere in the document (e.g.
# "Object-Centric Slot Diffusion") is left untouched.

$d = $word.ActiveDocument

# wdNoHighlight = 0 ; wdFindContinue = 1
$wdNoHighlight  = 0
$wdFindContinue = 1

$targetText = "SlotDiffusion: Object-Centric Generative Modeling with Diffusion Models"

# Preferred path: find the whole phrase (it spans two adjacent runs in the
# original markup, but Find treats the paragraph text as one contiguous
# string) and clear the highlight on that single range - this mirrors what
# a user does by selecting the title and clicking "No Color" on the text
# highlight button once.
$rng = $d.Content
$found = $rng.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if ($found -and $rng.Find.Found) {
    $rng.HighlightColorIndex = $wdNoHighlight
}

# Defensive fallback: in case the engine's Find could not match the full
# phrase as one span (e.g. due to the intervening spell-check markers),
# clear the highlight on each of the two runs' text individually too.
# These are no-ops if the text is already un-highlighted above.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("SlotDiffusion", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if ($found1 -and $rng1.Find.Found) {
    $rng1.HighlightColorIndex = $wdNoHighlight
}

$rng2 = $d.Content
$found2 = $rng2.Find.Execute(": Object-Centric Generative Modeling with Diffusion Models", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if ($found2 -and $rng2.Find.Found) {
    $rng2.HighlightColorIndex = $wdNoHighlight
}
